# Assignment5_Workbook.xlsx — "germanized assignments, minor typo fixes"
#
# 1. Translate the UI labels / headers on both sheets (and the matching
#    Excel Table column headers, which are kept in sync automatically by
#    Cells.Replace) from English to German.
# 2. Fix the "Cell Phone" header typo while we're at it (becomes lowercase
#    "telefon" per the source diff).
# 3. Change the Time column's number format from 12-hour "h:mm AM/PM" to a
#    24-hour "h:mm" display.
# 4. Rename the two worksheet tabs to their German equivalents.
# 5. Move the active selection/tab: the Player Info ("Spieler Info") sheet
#    becomes the active tab with K10 selected; Schedule ("Spielplan") is
#    left with D15 selected.

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$playerInfo = $wb.Worksheets.Item("Player Info")

$whole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

# --- Text / label translations -------------------------------------------
# Applies to literal cell text throughout a sheet (also keeps the
# corresponding Excel Table column names and the big banner title in sync).
$replacements = @{
    "Day"          = "Tag"
    "Date"         = "Datum"
    "Time"         = "Uhrzeit"
    "Opponent"     = "Gegner"
    "Saturday"     = "Samstag"
    "Friday"       = "Freitag"
    "Sunday"       = "Sonntag"
    "Cell Phone"   = "telefon"
    "Position(s)"  = "Position(en)"
    "First Name"   = "Vorname"
    "Last"         = "Nachname"
    "Street Address" = "Adresse"
    "Bulls Team Roster: Co-ed Softball 2013" = "Bulls Mannschaft: Softball 2013"
}

foreach ($sheet in @($schedule, $playerInfo)) {
    foreach ($key in $replacements.Keys) {
        $sheet.Cells.Replace($key, $replacements[$key], $whole)
    }
}

# --- Time column number format: 12h AM/PM -> 24h --------------------------
$schedule.Range("E1:E99").NumberFormat = "h:mm;@"

# --- Rename the worksheet tabs to German ----------------------------------
$schedule.Name = "Spielplan"
$playerInfo.Name = "Spieler Info"

# --- Update selections on each sheet, then make "Spieler Info" the active tab
$schedule.Range("D15").Select()
$playerInfo.Range("K10").Select()
$playerInfo.Activate()
